$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: normalize the font used by the numeric "Vars Used/# cPCA/Alpha"
#     cells so they match the rest of the table (explicit black text instead
#     of the theme-color font).
$ws.Range("K18").Font.Color = 0
$ws.Range("N18").Font.Color = 0
$ws.Range("O18").Font.Color = 0

# --- Row 19: this was a quote-prefixed placeholder row (every cell held an
#     empty, text-typed value). Fill it in with the results of running the
#     full dataset restricted to the BP-only variables.
$ws.Range("A19").Value = "ukb51139_subset.csv"
$ws.Range("B19").Value = "28012 x 4"
$ws.Range("C19").Value = "Sex/Age/BP only"
$ws.Range("D19").Value = "no events"
$ws.Range("E19").Value = "> 140/80"
$ws.Range("F19").Value = "zscore"
$ws.Range("G19").Value = "median"
$ws.Range("H19").Value = "none"
$ws.Range("I19").Value = 50
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = "47.7 & 35.8"
$ws.Range("M19").Value = "19.9 & 43.8"
$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 1.21

# --- Row 20: used to be the next quote-prefixed placeholder row; reset it to
#     a plain blank row (matching the formatting already used by the other
#     blank rows further down, e.g. row 24) now that row 19 is filled in.
$ws.Range("A24:O24").Copy()
$ws.Range("A20:O20").PasteSpecial(-4122)
$ws.Range("A20:H20").ClearContents()
$ws.Range("L20:M20").ClearContents()

# --- Remaining blank rows (21-23): same font normalization as row 18 so all
#     blank rows share one consistent (non theme-color) numeric style.
$ws.Range("I21:I23").Font.Color = 0
$ws.Range("K21:K23").Font.Color = 0
$ws.Range("N21:N23").Font.Color = 0
$ws.Range("O21:O23").Font.Color = 0

# --- Row-height normalization: every row in the table ends up at the same
#     19.5pt height.
$ws.Rows.Item(8).RowHeight = 19.5
$ws.Rows.Item(9).RowHeight = 19.5
$ws.Rows.Item(20).RowHeight = 19.5
$ws.Rows.Item(21).RowHeight = 19.5
$ws.Rows.Item(22).RowHeight = 19.5
$ws.Rows.Item(23).RowHeight = 19.5
